$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Column A: fix the "1.. Arrays" typo and add the new rows' numbering ---
$ws.Range("A1").Value = "1. Arrays"
$ws.Range("A2").Value = "2. Arrays"
$ws.Range("A3").Value = "3. Arrays"
$ws.Range("A4").Value = "4. Arrays"

# --- Column B: problem titles ---
$ws.Range("B2").Value = "Missing Number"
$ws.Range("B3").Value = "Desappeared numbers"
$ws.Range("B4").Value = "single Number"

# --- Column C: explanations (row 3 left blank) ---
$ws.Range("C2").Value = "using xor with the given array and numbers from 1 .. N in result the missing number will get. Because 1^1 = 0.
Another way is to calculate the sum of the array and the sum of the range and subtract."
$ws.Range("C4").Value = "using xor with the given array  1 .. N in result the missing number will get. Because 1^1 = 0.
Another way is to calculate the sum of the array and the 2* sum of the range and subtract."

# Reuse the existing highlighted-fill format from B1 for the rest of column B
$ws.Range("B1").Copy()
$ws.Range("B2:B4").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# --- Alignment / wrap formatting, applied in the order that recreates the
# expected style table (vertical-top-only, then vertical-top+wrap, then wrap-only) ---
$ws.Range("C3").VerticalAlignment = -4160
$ws.Range("C2").VerticalAlignment = -4160
$ws.Range("C2").WrapText = $true
$ws.Range("C4").WrapText = $true

# --- Row heights to match the final layout ---
$ws.Rows(2).RowHeight = 31.5
$ws.Rows(3).RowHeight = 27
$ws.Rows(4).RowHeight = 30

# --- Selection moves to C3 ---
$ws.Range("C3").Select()

Write-Output "done"
